# Azerbaijan Premier League workbook update (24-02-2024 12:40)
# 1) Rows 82-84 get their betting-odds data "rotated" (row82<-old83, row83<-old84, row84<-old82)
# 2) A new fixture row (205) is appended for a future match (no result yet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 82 (was id=80, match 5573342) now receives the former row-83 data ---
$ws.Range("B82").Value = 5573343
$ws.Range("F82").Value = "Shamakhi FK"
$ws.Range("G82").Value = "FK Gabala"
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = "D"
$ws.Range("K82").Value = 3.5
$ws.Range("L82").Value = 3.1
$ws.Range("M82").Value = 2
$ws.Range("N82").Value = 3.3
$ws.Range("O82").Value = 3.2
$ws.Range("P82").Value = 2.05
$ws.Range("Q82").Value = 0.25
$ws.Range("R82").Value = 2
$ws.Range("S82").Value = 1.8
$ws.Range("T82").Value = 2.5
$ws.Range("U82").Value = 1.975
$ws.Range("V82").Value = 1.825
$ws.Range("W82").Value = -1
$ws.Range("X82").Value = 2.2
$ws.Range("Y82").Value = -1
$ws.Range("Z82").Value = 0.5
$ws.Range("AA82").Value = -0.5
$ws.Range("AB82").Value = -1
$ws.Range("AC82").Value = 0.825

# --- Row 83 (was id=81, match 5573343) now receives the former row-84 data ---
$ws.Range("B83").Value = 5574442
$ws.Range("F83").Value = "FK Qarabag"
$ws.Range("G83").Value = "FK Sumqayit"
$ws.Range("H83").Value = 1
$ws.Range("I83").Value = 2
$ws.Range("J83").Value = "A"
$ws.Range("K83").Value = 1.125
$ws.Range("L83").Value = 7.5
$ws.Range("M83").Value = 15
$ws.Range("N83").Value = 1.2
$ws.Range("O83").Value = 6
$ws.Range("P83").Value = 11
$ws.Range("Q83").Value = -2.25
$ws.Range("R83").Value = 1.975
$ws.Range("S83").Value = 1.825
$ws.Range("T83").Value = 3.5
$ws.Range("U83").Value = 1.825
$ws.Range("V83").Value = 1.975
$ws.Range("W83").Value = -1
$ws.Range("X83").Value = -1
$ws.Range("Y83").Value = 10
$ws.Range("Z83").Value = -1
$ws.Range("AA83").Value = 0.825
$ws.Range("AB83").Value = -1
$ws.Range("AC83").Value = 0.9750000000000001

# --- Row 84 (was id=82, match 5574442) now receives the former row-82 data ---
$ws.Range("B84").Value = 5573342
$ws.Range("F84").Value = "PFK Turan Tovuz"
$ws.Range("G84").Value = "Sabail FC"
$ws.Range("H84").Value = 2
$ws.Range("I84").Value = 2
$ws.Range("J84").Value = "D"
$ws.Range("K84").Value = 2.6
$ws.Range("L84").Value = 3
$ws.Range("M84").Value = 2.6
$ws.Range("N84").Value = 2.8
$ws.Range("O84").Value = 2.875
$ws.Range("P84").Value = 2.5
$ws.Range("Q84").Value = 0
$ws.Range("R84").Value = 2.05
$ws.Range("S84").Value = 1.75
$ws.Range("T84").Value = 2.25
$ws.Range("U84").Value = 1.875
$ws.Range("V84").Value = 1.925
$ws.Range("W84").Value = -1
$ws.Range("X84").Value = 1.875
$ws.Range("Y84").Value = -1
$ws.Range("Z84").Value = 0
$ws.Range("AA84").Value = 0
$ws.Range("AB84").Value = 0.875
$ws.Range("AC84").Value = -1

# --- New row 205: upcoming fixture FK Qarabag vs FK Sumqayit (no score/result yet) ---
# Copy formatting (styles/number-format) from the last existing data row, then
# clear it and fill in only the columns that the new fixture actually has.
$ws.Range("A204:AC204").Copy($ws.Range("A205:AC205"))
$ws.Range("B205:AC205").ClearContents()

$ws.Range("A205").Value = 203
$ws.Range("B205").Value = 7809824
$ws.Range("C205").Value = "Azerbaijan Premier League"
$ws.Range("D205").Value = "Azerbaijan Premier League"
$ws.Range("E205").Value = 45347.5
$ws.Range("F205").Value = "FK Qarabag"
$ws.Range("G205").Value = "FK Sumqayit"
$ws.Range("K205").Value = 1.333
$ws.Range("L205").Value = 4.333
$ws.Range("M205").Value = 7.5
$ws.Range("N205").Value = 1.25
$ws.Range("O205").Value = 4.75
$ws.Range("P205").Value = 10
$ws.Range("Q205").Value = -1.75
$ws.Range("R205").Value = 1.875
$ws.Range("S205").Value = 1.925
$ws.Range("T205").Value = 3
$ws.Range("U205").Value = 1.875
$ws.Range("V205").Value = 1.925
$ws.Range("W205").Value = 0
$ws.Range("X205").Value = 0
$ws.Range("Y205").Value = 0
$ws.Range("Z205").Value = 0
$ws.Range("AA205").Value = 0
